# Insert a new row above the current row 9 ("Panos" / ...) so that it
# becomes row 10 and everything below shifts down by one. The freshly
# inserted row 8 gets a single value in column A: "preconInput".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").EntireRow.Insert()
$ws.Range("A8").Value = "preconInput"

# Move/record the active selection on the newly added cell, matching the
# updated <selection activeCell="A8" ... sqref="A8"/> in the sheet view.
[void]$ws.Range("A8").Select()
